$wb = $excel.ActiveWorkbook

# A worksheet already carrying the "m/d/yyyy" date style (cellXfs index 2) so
# we can reuse that style (instead of minting a brand-new numFmt) when we
# stamp the "event_day" column on the new "Day 21" sheet.
$dateStyleSource = $wb.Worksheets.Item("Day 4")

# Both new tabs get appended after the current last sheet ("Day 16"), in
# order, so sheetId/r:id allocation and tab order line up with the target.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$day19 = $wb.Worksheets.Add($null, $lastSheet)
$day19.Name = "Day 19"

$day21 = $wb.Worksheets.Add($null, $day19)
$day21.Name = "Day 21"

# ---------------------------------------------------------------------
# Day 19: account_id / income
# ---------------------------------------------------------------------
$day19.Range("A1").Value = "account_id"
$day19.Range("B1").Value = "income"

$day19Rows = @(
    @(3, 108939),
    @(2, 12747),
    @(8, 87709),
    @(6, 91796)
)
for ($i = 0; $i -lt $day19Rows.Count; $i++) {
    $r = $i + 2
    $day19.Cells.Item($r, 1).Value = $day19Rows[$i][0]
    $day19.Cells.Item($r, 2).Value = $day19Rows[$i][1]
}

# ---------------------------------------------------------------------
# Day 21: emp_id / event_day / in_time / out_time
# Columns are populated A, C, D, then B so the shared-string table picks
# up "in_time"/"out_time" before "event_day", matching the source order.
# ---------------------------------------------------------------------
$day21.Range("A1").Value = "emp_id"
$day21.Range("C1").Value = "in_time"
$day21.Range("D1").Value = "out_time"
$day21.Range("B1").Value = "event_day"

$day21Rows = @(
    @(1, "2020-11-28", 4, 32),
    @(1, "2020-11-28", 55, 200),
    @(1, "2020-12-03", 1, 42),
    @(2, "2020-11-28", 3, 33),
    @(2, "2020-12-09", 47, 74)
)
for ($i = 0; $i -lt $day21Rows.Count; $i++) {
    $r = $i + 2
    $day21.Cells.Item($r, 1).Value = $day21Rows[$i][0]
    $day21.Cells.Item($r, 2).Value = $day21Rows[$i][1]
    $day21.Cells.Item($r, 3).Value = $day21Rows[$i][2]
    $day21.Cells.Item($r, 4).Value = $day21Rows[$i][3]
}

# Apply the existing date number format to the event_day column by copying
# it from a cell that already carries it, so the style is reused rather
# than duplicated.
$dateStyleSource.Range("D2").Copy() | Out-Null
$day21.Range("B2:B6").PasteSpecial(-4122) | Out-Null

# Match the saved selections / active tab from the target workbook.
$day21.Range("B2").Select() | Out-Null
$day19.Activate() | Out-Null
$day19.Range("G9").Select() | Out-Null
